$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 14.23347321995233
$ws.Range("R2").Value = 128.101258979571
$ws.Range("S2").Value = 0.6311762527593259
$ws.Range("T2").Value = 0.6311762527593258

# Row 3
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 2.013415967717334
$ws.Range("S3").Value = 0.08928392431779728
$ws.Range("T3").Value = 0.08928392431779726

# Row 4
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 6.303821739313999
$ws.Range("R4").Value = 56.73439565382599
$ws.Range("S4").Value = 0.2795398229228769
$ws.Range("T4").Value = 0.2795398229228769
